$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1689724
$ws.Range("J17").Value = 1689724
$ws.Range("L17").Value = 5069172
$ws.Range("N17").Value = -5069508

$ws.Range("H88").Value = 8016.5625
$ws.Range("J88").Value = 8840.357
$ws.Range("L88").Value = 8840.357
$ws.Range("N88").Value = -9652.357

$ws.Range("H91").Value = 8016.5625
$ws.Range("J91").Value = 8840.357
$ws.Range("L91").Value = 8840.357
$ws.Range("N91").Value = -11648.357

$ws.Range("H107").Value = 11364482
$ws.Range("I107").Value = 11905171
$ws.Range("J107").Value = 10000
$ws.Range("K107").Value = 11905171
$ws.Range("L107").Value = 10000
$ws.Range("M107").Value = -11903251
$ws.Range("N107").Value = -13840

$ws.Range("H116").Value = 7309.737
$ws.Range("I116").Value = 10990.909
$ws.Range("J116").Value = 2248.125
$ws.Range("K116").Value = 10990.909
$ws.Range("L116").Value = 2248.125
$ws.Range("M116").Value = -7548.909
$ws.Range("N116").Value = -9132.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 746.3570999999999
$ws.Range("I2").Value = 469.66666
$ws.Range("K2").Value = 469.66666
$ws.Range("M2").Value = -356.66666

$ws.Range("H8").Value = 4268.6665
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 4268.6665
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 4268.6665
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = -4556.6665

$ws.Range("H61").Value = 2849.825
$ws.Range("I61").Value = 2027.5769
$ws.Range("J61").Value = 4376.857
$ws.Range("K61").Value = 2027.5769
$ws.Range("L61").Value = 4376.857
$ws.Range("M61").Value = -1815.5769
$ws.Range("N61").Value = -4800.857

$ws.Range("H110").Value = 1477.9565
$ws.Range("I110").Value = 1257.6666
$ws.Range("K110").Value = 1257.6666
$ws.Range("M110").Value = 787.3334

$ws.Range("H116").Value = 746.3570999999999
$ws.Range("I116").Value = 469.66666
$ws.Range("K116").Value = 469.66666
$ws.Range("M116").Value = 1824.33334

$ws.Range("H119").Value = 55739.6
$ws.Range("J119").Value = 55739.6
$ws.Range("L119").Value = 55739.6
$ws.Range("N119").Value = -65415.6

$ws.Range("H122").Value = 953132.3
$ws.Range("I122").Value = 1224522.8
$ws.Range("J122").Value = 3265.6667
$ws.Range("K122").Value = 3673568.4
$ws.Range("L122").Value = 9797.000100000001
$ws.Range("M122").Value = -3671118.4
$ws.Range("N122").Value = -14697.0001

$ws.Range("H132").Value = 3128209.8
$ws.Range("I132").Value = 1816.9
$ws.Range("J132").Value = 8338864.5
$ws.Range("K132").Value = 5450.700000000001
$ws.Range("L132").Value = 25016593.5
$ws.Range("M132").Value = -2920.700000000001
$ws.Range("N132").Value = -25021653.5

$ws.Range("H136").Value = 2849.825
$ws.Range("I136").Value = 2027.5769
$ws.Range("J136").Value = 4376.857
$ws.Range("K136").Value = 6082.7307
$ws.Range("L136").Value = 13130.571
$ws.Range("M136").Value = -3532.7307
$ws.Range("N136").Value = -18230.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 746.3570999999999
$ws.Range("I3").Value = 469.66666
$ws.Range("K3").Value = 469.66666
$ws.Range("M3").Value = -355.66666

$ws.Range("H107").Value = 1160.1852
$ws.Range("I107").Value = 1078.2142
$ws.Range("K107").Value = 1078.2142
$ws.Range("M107").Value = 841.7858000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2264338
$ws.Range("J16").Value = 2219.6428
$ws.Range("L16").Value = 2219.6428
$ws.Range("N16").Value = -2793.6428

$ws.Range("H107").Value = 1571.4286
$ws.Range("I107").Value = 433.33334
$ws.Range("J107").Value = 2425
$ws.Range("K107").Value = 433.33334
$ws.Range("L107").Value = 2425
$ws.Range("M107").Value = 1486.66666
$ws.Range("N107").Value = -6265

$ws.Range("H113").Value = 2264338
$ws.Range("J113").Value = 2219.6428
$ws.Range("L113").Value = 2219.6428
$ws.Range("N113").Value = -6559.6428

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 3575124.5
$ws.Range("I33").Value = 4762030
$ws.Range("J33").Value = 14407.857
$ws.Range("K33").Value = 28572180
$ws.Range("L33").Value = 86447.14199999999
$ws.Range("M33").Value = -28571897
$ws.Range("N33").Value = -87013.14199999999

$ws.Range("H44").Value = 766.86664
$ws.Range("I44").Value = 343.2857
$ws.Range("J44").Value = 1137.5
$ws.Range("K44").Value = 1029.8571
$ws.Range("L44").Value = 3412.5
$ws.Range("M44").Value = -631.8571000000002
$ws.Range("N44").Value = -4208.5

$ws.Range("H80").Value = 2855
$ws.Range("J80").Value = 2855
$ws.Range("L80").Value = 8565
$ws.Range("N80").Value = -10437

$ws.Range("H83").Value = 2855
$ws.Range("J83").Value = 2855
$ws.Range("L83").Value = 25695
$ws.Range("N83").Value = -35055

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 33660
$ws.Range("I34").Value = 25980
$ws.Range("J34").Value = 37500
$ws.Range("K34").Value = 25980
$ws.Range("L34").Value = 37500
$ws.Range("M34").Value = -25712
$ws.Range("N34").Value = -38036

$ws.Range("H76").Value = 33660
$ws.Range("I76").Value = 25980
$ws.Range("J76").Value = 37500
$ws.Range("K76").Value = 25980
$ws.Range("L76").Value = 37500
$ws.Range("M76").Value = -25665
$ws.Range("N76").Value = -38130

$ws.Range("H79").Value = 33660
$ws.Range("I79").Value = 25980
$ws.Range("J79").Value = 37500
$ws.Range("K79").Value = 25980
$ws.Range("L79").Value = 37500
$ws.Range("M79").Value = -24888
$ws.Range("N79").Value = -39684

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 13613699
$ws.Range("I82").Value = 2501247.5
$ws.Range("K82").Value = 2501247.5
$ws.Range("M82").Value = -2500886.5

$ws.Range("H85").Value = 13613699
$ws.Range("I85").Value = 2501247.5
$ws.Range("K85").Value = 2501247.5
$ws.Range("M85").Value = -2499999.5

$ws.Range("H119").Value = 30000
$ws.Range("J119").Value = 30000
$ws.Range("L119").Value = 30000
$ws.Range("N119").Value = -39676

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1061.7142
$ws.Range("J81").Value = 1500
$ws.Range("L81").Value = 3000
$ws.Range("N81").Value = -5122

$ws.Range("H84").Value = 1061.7142
$ws.Range("J84").Value = 1500
$ws.Range("L84").Value = 15000
$ws.Range("N84").Value = -25608

$ws.Range("H119").Value = 66000
$ws.Range("J119").Value = 66000
$ws.Range("L119").Value = 66000
$ws.Range("N119").Value = -75676

$ws.Range("H136").Value = 2986.4595
$ws.Range("I136").Value = 3932
$ws.Range("K136").Value = 11796
$ws.Range("M136").Value = -9246
